# A new weekly price record (Choclero, "Primera", Región del Maule) was
# added to the daily log for Macroferia Regional de Talca - Choclo.
# In the source table this shows up as a brand-new row inserted right
# above the existing row 91, pushing every following record down by one
# row (old row 91 -> new row 92, ..., old row 209 -> new row 210).
#
# The newly inserted row duplicates the record that used to sit at row 91
# (same market/region/variety/quality/prices/unit/origin/classification)
# except for the reporting date (column D) and the traded volume
# (column J), which carry the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 91; everything from old row 91
# downward shifts to row+1 (old row 91 is now row 92, etc.), and column
# formatting (e.g. the date style on column D) carries down with it.
$ws.Rows("91:91").Insert()

# Seed the new row 91 with a copy of what is now row 92 (the record that
# used to be row 91 before the insert), column by column.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(91, $col).Value = $ws.Cells.Item(92, $col).Value2
}

# Overwrite the two fields that actually differ for the new record:
#   D91 -> Fecha (date serial 44629 = 2022-03-09)
#   J91 -> Volumen (40000)
$ws.Cells.Item(91, 4).Value = 44629
$ws.Cells.Item(91, 10).Value = 40000
